# Country-Code.xlsx update:
# Adds "Currency Code" (C) and "Conversion Rate" (D) columns with the
# per-country currency symbol and a USD conversion rate, plus light
# banding borders + column widths matching the authored worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ------------------------------------------------------
$ws.Range("C1").Value = "Currency Code"
$ws.Range("D1").Value = "Conversion Rate"

# ---- Currency code + conversion-rate data (rows 2-16) ----------------
$ws.Range("C2").Value = "INR"
$ws.Range("D2").Value = 0.01287823000000000

$ws.Range("C3").Value = "AUD"
$ws.Range("D3").Value = 0.70033959000000001

$ws.Range("C4").Value = "BRL"
$ws.Range("D4").Value = 0.20153618000000001

$ws.Range("C5").Value = "CAD"
$ws.Range("D5").Value = 0.77912552000000002

$ws.Range("C6").Value = "IDR"
$ws.Range("D6").Value = 0.00006787000000000

$ws.Range("C7").Value = "NZD"
$ws.Range("D7").Value = 0.63404706000000000

$ws.Range("C8").Value = "PHP"
$ws.Range("D8").Value = 0.01910187000000000

$ws.Range("C9").Value = "QAR"
$ws.Range("D9").Value = 0.27368767999999999

$ws.Range("C10").Value = "SGD"
$ws.Range("D10").Value = 0.72051524000000000

$ws.Range("C11").Value = "ZAR"
$ws.Range("D11").Value = 0.06264032000000000

$ws.Range("C12").Value = "LKR"
$ws.Range("D12").Value = 0.00277712000000000

$ws.Range("C13").Value = "TRY"
$ws.Range("D13").Value = 0.06269716000000000

$ws.Range("C14").Value = "AED"
$ws.Range("D14").Value = 0.27224554000000001

$ws.Range("C15").Value = "GBP"
$ws.Range("D15").Value = 1.24023062999999989

$ws.Range("C16").Value = "USD"
$ws.Range("D16").Value = 1.00000000000000000

# ---- Banding borders on the conversion-rate column (D2:D15) ----------
# Thin top+bottom rule in the workbook accent color (Accent 6, ~40% tint,
# i.e. theme color 9 / tint 0.4 ~= RGB A9D18E)
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 4)

    $top = $cell.Borders.Item(8)
    $top.LineStyle = 1
    $top.Color = 9359785

    $bottom = $cell.Borders.Item(9)
    $bottom.LineStyle = 1
    $bottom.Color = 9359785
}

# ---- Column widths -----------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.333333333333332
$ws.Columns.Item(4).ColumnWidth = 13.833333333333332
